$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Write a literal text string into a cell without Excel's autodetect turning
# it into a date/number (e.g. "1888-08-06") and without minting a brand new
# cell style in the process: assign it as a quoted-string formula first (the
# formula result inherits the cell's normal, column-driven style), then
# Copy + PasteSpecial(xlPasteValues) it onto itself to freeze the formula
# down to a plain shared-string value.
function Set-TextLiteral($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# Row 86
$ws.Cells.Item(86, 2).Value = 8
Set-TextLiteral $ws.Cells.Item(86, 3) "1888-08-06"
Set-TextLiteral $ws.Cells.Item(86, 4) "Other Article"
Set-TextLiteral $ws.Cells.Item(86, 5) "N"
Set-TextLiteral $ws.Cells.Item(86, 6) "Told in the Pulpit."

# Row 87
$ws.Cells.Item(87, 2).Value = 1
Set-TextLiteral $ws.Cells.Item(87, 3) "1888-08-07"
Set-TextLiteral $ws.Cells.Item(87, 4) "CSG Article"
Set-TextLiteral $ws.Cells.Item(87, 5) "N"

# Row 88
$ws.Cells.Item(88, 2).Value = 2
Set-TextLiteral $ws.Cells.Item(88, 3) "1888-08-07"
Set-TextLiteral $ws.Cells.Item(88, 4) "CSG Article"
Set-TextLiteral $ws.Cells.Item(88, 5) "N"

$ws.Range("F88").Select() | Out-Null
